$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1170
$ws1.Range("F4").Value = 2631
$ws1.Range("F5").Value = 228

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1170
$ws4.Range("F6").Value = 2631
$ws4.Range("F8").Value = 228
